# Update the Status for "bb20a7ba-c30e-467a-8914-ca89bacca223.md" from
# "Ready for handoff" to "In Translation" across the Overview, zh-cn and
# de-de sheets (row 4 in each sheet).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B4").Value = "In Translation"
$overview.Range("C4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B4").Value = "In Translation"
